# Dictionary.xlsx update — refresh the "Programming Skills" list on Sheet1.
#
# The header in F2 is renamed from "Programming Skills" to "skills", the
# language-specific entries "r", "c#" and "C++" are dropped from the list
# (the remaining entries shift up to fill the gap) and several new
# technologies used in the R-code project are appended at the bottom of
# the column: spark, excel (again), aws, hadoop, azure, rstudio, hive.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column F ("Programming Skills" -> "skills") ---------------------------
$ws.Range("F2").Value  = "skills"

$ws.Range("F8").Value  = "tableau"
$ws.Range("F9").Value  = "django"
$ws.Range("F10").Value = "javascript"
$ws.Range("F11").Value = "excel"
$ws.Range("F12").Value = "sas"
$ws.Range("F13").Value = "matlab"
$ws.Range("F14").Value = "github"
$ws.Range("F15").Value = "spark"
$ws.Range("F16").Value = "excel"
$ws.Range("F17").Value = "aws"
$ws.Range("F18").Value = "hadoop"
$ws.Range("F19").Value = "azure"
$ws.Range("F20").Value = "rstudio"
$ws.Range("F21").Value = "hive"

# --- Cosmetic view-state refresh (matches the author's last selection) ----
$ws.Activate()
$ws.Range("G20").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

# --- Theme rename (author's Office install uses a localized theme name) ---
$wb.Theme.Name = "Thème Office"
